$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary block with live formulas (rows 14-17) ---
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Bold/size-12/vertically-centred style for the new B14:B17 summary values.
foreach ($r in 14..17) {
    $cell = $ws.Cells.Item($r, 2)
    $f = $cell.Font
    $f.Bold = $true
    $f.Size = 12
    $cell.VerticalAlignment = -4108
}
$ws.Rows("14:17").RowHeight = 15.6

# --- Secondary summary block with hard-coded values (rows 20-23) ---
$ws.Range("A20").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B20").Value = 0.996212121

$ws.Range("A21").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B21").Value = 1.000865052

$ws.Range("A22").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B22").Value = 0.962121212

$ws.Range("A23").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B23").Value = 1.008650519

# --- Page setup / selection cosmetics ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("B17").Select()
